$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.310603618621826
$ws.Range("B1").Value = 2.314539194107056
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.766597390174866
$ws.Range("E1").Value = 1.011873126029968
